$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,12
$row2[0,0] = 0.00244238008256925
$row2[0,1] = 0.00353318336545636
$row2[0,2] = 0.00617524990083813
$row2[0,3] = 0.00951957580678674
$row2[0,4] = 0.0116637630172437
$row2[0,5] = 0.0132099052977561
$row2[0,6] = 0.0141206308298382
$row2[0,7] = 0.0146192594223433
$row2[0,8] = 0.0147229778038046
$row2[0,9] = 0.0144074050157856
$row2[0,10] = 0.0140066051218847
$row2[0,11] = 0.013451159753295
$ws.Range("A2:L2").Value = $row2

$row3 = New-Object "object[,]" 1,12
$row3[0,0] = 0.000151465284429442
$row3[0,1] = 0.00250150911228125
$row3[0,2] = 0.00561083144496203
$row3[0,3] = 0.00728626696092478
$row3[0,4] = 0.00850976591351137
$row3[0,5] = 0.00917284549710259
$row3[0,6] = 0.00949232425422627
$row3[0,7] = 0.0094585004449071
$row3[0,8] = 0.00907066273920572
$row3[0,9] = 0.00865493781180123
$row3[0,10] = 0.00828671601764613
$row3[0,11] = 0.00800684040539673
$ws.Range("A3:L3").Value = $row3

$row4 = New-Object "object[,]" 1,12
$row4[0,0] = 0.0038145431900789
$row4[0,1] = 0.00718610268062892
$row4[0,2] = 0.00876378455720842
$row4[0,3] = 0.00988601526406115
$row4[0,4] = 0.0104324394395268
$row4[0,5] = 0.0106429711238965
$row4[0,6] = 0.010503406304823
$row4[0,7] = 0.0100069526093559
$row4[0,8] = 0.00949920118741533
$row4[0,9] = 0.00903673469851177
$row4[0,10] = 0.00867524891595853
$row4[0,11] = 0.00834654286506908
$ws.Range("A4:L4").Value = $row4

$row5 = New-Object "object[,]" 1,12
$row5[0,0] = 0.00238972992031528
$row5[0,1] = 0.00221310583839735
$row5[0,2] = 0.00221012176776939
$row5[0,3] = 0.00205532699523641
$row5[0,4] = 0.00187561465298415
$row5[0,5] = 0.00175038464538245
$row5[0,6] = 0.0024296298810152
$row5[0,7] = 0.00324765085894375
$row5[0,8] = 0.0049892327363201
$row5[0,9] = 0.00625817632601084
$row5[0,10] = 0.0070444854680022
$row5[0,11] = 0.00713962429770822
$ws.Range("A5:L5").Value = $row5

$row6 = New-Object "object[,]" 1,12
$row6[0,0] = 0.00239186719864204
$row6[0,1] = 0.00305884993791666
$row6[0,2] = 0.00385554646530816
$row6[0,3] = 0.00449957193528507
$row6[0,4] = 0.00530016577781364
$row6[0,5] = 0.00672297185959227
$row6[0,6] = 0.00792748970224845
$row6[0,7] = 0.009786373639006
$row6[0,8] = 0.0112055326058074
$row6[0,9] = 0.012135033283672
$row6[0,10] = 0.0123361816496052
$row6[0,11] = 0.0120313081845058
$ws.Range("A6:L6").Value = $row6

$row7 = New-Object "object[,]" 1,12
$row7[0,0] = 0.00082704782788845
$row7[0,1] = 0.000787976574078821
$row7[0,2] = 0.000723282601183167
$row7[0,3] = 0.000750728845763598
$row7[0,4] = 0.00216088287106792
$row7[0,5] = 0.00325311461435729
$row7[0,6] = 0.00532802183089969
$row7[0,7] = 0.00673756584915684
$row7[0,8] = 0.00757494075199025
$row7[0,9] = 0.00763072551074612
$row7[0,10] = 0.00729747765576614
$row7[0,11] = 0.00699529471235806
$ws.Range("A7:L7").Value = $row7

$row8 = New-Object "object[,]" 1,12
$row8[0,0] = 0.000782994922419178
$row8[0,1] = 0.00116087654234229
$row8[0,2] = 0.00202970682107659
$row8[0,3] = 0.00390184249357477
$row8[0,4] = 0.00522500241994879
$row8[0,5] = 0.00741973461106876
$row8[0,6] = 0.00891669081835738
$row8[0,7] = 0.00979674009321217
$row8[0,8] = 0.00984572302807704
$row8[0,9] = 0.00942430589746803
$row8[0,10] = 0.00904019173154178
$row8[0,11] = 0.00866663731624204
$ws.Range("A8:L8").Value = $row8

$row9 = New-Object "object[,]" 1,12
$row9[0,0] = 0.00000181765213103091
$row9[0,1] = 0.000855522332326968
$row9[0,2] = 0.00286618734138847
$row9[0,3] = 0.004069898388691
$row9[0,4] = 0.00637054561448776
$row9[0,5] = 0.00783763200089958
$row9[0,6] = 0.00864496279445322
$row9[0,7] = 0.00858912803858337
$row9[0,8] = 0.0081241850500947
$row9[0,9] = 0.00771815872837471
$row9[0,10] = 0.0073602762116692
$row9[0,11] = 0.00705185123861307
$ws.Range("A9:L9").Value = $row9

$row10 = New-Object "object[,]" 1,12
$row10[0,0] = 0.00121331059970853
$row10[0,1] = 0.00351404820829779
$row10[0,2] = 0.00470373851669452
$row10[0,3] = 0.00712686633124086
$row10[0,4] = 0.00859030229687423
$row10[0,5] = 0.00934239473762207
$row10[0,6] = 0.00918701159605598
$row10[0,7] = 0.00862166676930996
$row10[0,8] = 0.00814013648782111
$row10[0,9] = 0.007723746298979
$row10[0,10] = 0.007369495741074
$row10[0,11] = 0.00715909802130062
$ws.Range("A10:L10").Value = $row10

$row11 = New-Object "object[,]" 1,12
$row11[0,0] = 0.00257932366318858
$row11[0,1] = 0.00307889478584938
$row11[0,2] = 0.00544244263696173
$row11[0,3] = 0.00665825963247505
$row11[0,4] = 0.00718635654103573
$row11[0,5] = 0.00685180430447383
$row11[0,6] = 0.00636166999040542
$row11[0,7] = 0.00598623984835467
$row11[0,8] = 0.00576544701070538
$row11[0,9] = 0.00551103549360673
$row11[0,10] = 0.00526597920314074
$row11[0,11] = 0.00510193660461917
$ws.Range("A11:L11").Value = $row11

$row12 = New-Object "object[,]" 1,12
$row12[0,0] = 0.00126029812416029
$row12[0,1] = 0.00172749904744152
$row12[0,2] = 0.00194508471229289
$row12[0,3] = 0.00186823980848923
$row12[0,4] = 0.0020046558944836
$row12[0,5] = 0.0038826114074717
$row12[0,6] = 0.00487427880290783
$row12[0,7] = 0.00581418914389486
$row12[0,8] = 0.00619777856079093
$row12[0,9] = 0.0061163104141968
$row12[0,10] = 0.00594101381126985
$row12[0,11] = 0.00569566838362121
$ws.Range("A12:L12").Value = $row12

$row13 = New-Object "object[,]" 1,12
$row13[0,0] = 0.00442246525295964
$row13[0,1] = 0.00492683909237924
$row13[0,2] = 0.00497752347249412
$row13[0,3] = 0.00434727909661704
$row13[0,4] = 0.00444435630004261
$row13[0,5] = 0.0045867786509448
$row13[0,6] = 0.0049748897053612
$row13[0,7] = 0.00503393729930913
$row13[0,8] = 0.00479972171711227
$row13[0,9] = 0.00455753598788944
$row13[0,10] = 0.00439269715021233
$row13[0,11] = 0.00431206330180497
$ws.Range("A13:L13").Value = $row13

$row14 = New-Object "object[,]" 1,12
$row14[0,0] = 0.00278346689407627
$row14[0,1] = 0.00447020469294412
$row14[0,2] = 0.00734833461178183
$row14[0,3] = 0.0108118390781476
$row14[0,4] = 0.0125348914581251
$row14[0,5] = 0.0138419976555211
$row14[0,6] = 0.0143730456968338
$row14[0,7] = 0.0142769880371125
$row14[0,8] = 0.0140244289520245
$row14[0,9] = 0.0136016006799759
$row14[0,10] = 0.0131773403734896
$row14[0,11] = 0.0127137354384076
$ws.Range("A14:L14").Value = $row14

$row15 = New-Object "object[,]" 1,11
$row15[0,0] = 0.000541070699749291
$row15[0,1] = 0.00305956467141557
$row15[0,2] = 0.00620093274004149
$row15[0,3] = 0.00728876782021925
$row15[0,4] = 0.00820497455123487
$row15[0,5] = 0.00843436012189725
$row15[0,6] = 0.00814533538109045
$row15[0,7] = 0.00778948649050051
$row15[0,8] = 0.00737002200901195
$row15[0,9] = 0.00699575882161242
$row15[0,10] = 0.00668658827652037
$ws.Range("A15:K15").Value = $row15

$row16 = New-Object "object[,]" 1,10
$row16[0,0] = 0.00329697652378425
$row16[0,1] = 0.00645958208420073
$row16[0,2] = 0.00715511307191236
$row16[0,3] = 0.00783761385620087
$row16[0,4] = 0.00786378699202496
$row16[0,5] = 0.00742495233998433
$row16[0,6] = 0.00697642336254542
$row16[0,7] = 0.00652795114963225
$row16[0,8] = 0.00615802683630032
$row16[0,9] = 0.00591332958338093
$ws.Range("A16:J16").Value = $row16

$row17 = New-Object "object[,]" 1,9
$row17[0,0] = 0.00246239583285046
$row17[0,1] = 0.00176646086741299
$row17[0,2] = 0.00148489847955028
$row17[0,3] = 0.00145698806224479
$row17[0,4] = 0.00244109081524043
$row17[0,5] = 0.00325877880802147
$row17[0,6] = 0.00431837606283988
$row17[0,7] = 0.00513843598471056
$row17[0,8] = 0.00609007191352512
$ws.Range("A17:I17").Value = $row17

$row18 = New-Object "object[,]" 1,8
$row18[0,0] = 0.00413016430164603
$row18[0,1] = 0.00480016863226276
$row18[0,2] = 0.00612794875125237
$row18[0,3] = 0.00788428239519262
$row18[0,4] = 0.00908109849908998
$row18[0,5] = 0.01028422890519
$row18[0,6] = 0.0111681518601699
$row18[0,7] = 0.0121073468236984
$ws.Range("A18:H18").Value = $row18

$row19 = New-Object "object[,]" 1,7
$row19[0,0] = 0.00222638612047099
$row19[0,1] = 0.00204439585242036
$row19[0,2] = 0.00168302346678072
$row19[0,3] = 0.00155499558239394
$row19[0,4] = 0.00204219524498947
$row19[0,5] = 0.00255444174360513
$row19[0,6] = 0.0034327560429624
$ws.Range("A19:G19").Value = $row19

$row20 = New-Object "object[,]" 1,6
$row20[0,0] = 0.00225246314586158
$row20[0,1] = 0.00436285801473914
$row20[0,2] = 0.00544908705158619
$row20[0,3] = 0.00674438494278162
$row20[0,4] = 0.00764740268009406
$row20[0,5] = 0.00867529168398467
$ws.Range("A20:F20").Value = $row20

$row21 = New-Object "object[,]" 1,5
$row21[0,0] = 0.00160739893459638
$row21[0,1] = 0.00166312995500645
$row21[0,2] = 0.00244167362106232
$row21[0,3] = 0.00296288977174058
$row21[0,4] = 0.00387114449110471
$ws.Range("A21:E21").Value = $row21

$row22 = New-Object "object[,]" 1,4
$row22[0,0] = 0.0012323470459803
$row22[0,1] = 0.00090563429440204
$row22[0,2] = 0.000748845719770071
$row22[0,3] = 0.00110794707037471
$ws.Range("A22:D22").Value = $row22

$row23 = New-Object "object[,]" 1,3
$row23[0,0] = 0.00191659447302683
$row23[0,1] = 0.00237823762034132
$row23[0,2] = 0.00354819192969732
$ws.Range("A23:C23").Value = $row23

$row24 = New-Object "object[,]" 1,2
$row24[0,0] = 0.00075348808863307
$row24[0,1] = 0.00065342594309793
$ws.Range("A24:B24").Value = $row24

$row25 = New-Object "object[,]" 1,1
$row25[0,0] = 0.00191779047916008
$ws.Range("A25:A25").Value = $row25
